$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Every existing row (2..470) got its "Förändrad" (column C) date bumped
#    from 45182 to 45184.
$ws.Range("C2:C470").Value2 = 45184

# 2. Row 470 gained an explicit row height (ht="15" customHeight="1"),
#    matching the style already used by the rows above it.
$ws.Rows.Item(470).RowHeight = 15

# 3. A brand-new row 471 was appended with a fresh record.
$ws.Range("A471").Value2 = "A 42809-2023"

$ws.Range("B471").Value2 = 45182
$ws.Range("B471").NumberFormat = $ws.Range("B470").NumberFormat

$ws.Range("C471").Value2 = 45184
$ws.Range("C471").NumberFormat = $ws.Range("C470").NumberFormat

$ws.Range("D471").Value2 = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E471").Value2 = "MARK"
$ws.Range("F471").Value2 = "Sveaskog"
$ws.Range("G471").Value2 = 7.6
$ws.Range("H471").Value2 = 0
$ws.Range("I471").Value2 = 0
$ws.Range("J471").Value2 = 0
$ws.Range("K471").Value2 = 0
$ws.Range("L471").Value2 = 0
$ws.Range("M471").Value2 = 0
$ws.Range("N471").Value2 = 0
$ws.Range("O471").Value2 = 0
$ws.Range("P471").Value2 = 0
$ws.Range("Q471").Value2 = 0

# R column uses the wrap-text style seen throughout the sheet but stays blank.
$ws.Range("R471").WrapText = $ws.Range("R470").WrapText
